# Updated cryptos list on Mon Jul  8 23:13:14 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) / Volume(1h) (E) refresh, plus a couple of rank swaps (rows 35/36, 49/50).
# Numeric-looking price strings are apostrophe-prefixed so Excel keeps them as text
# (matching the source data, which stores every cell as a string) instead of coercing
# them into real numbers and losing formatting (e.g. trailing zeros).

$ws.Range("D2").Value = "56.721.96"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "3.024.30"
$ws.Range("E3").Value = "  +2.74%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'510.60"
$ws.Range("E5").Value = "  +3.31%  "
$ws.Range("D6").Value = "'140.00"
$ws.Range("E6").Value = "  +4.81%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +2.14%  "
$ws.Range("D9").Value = "'7.13"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  +2.57%  "
$ws.Range("D11").Value = "'0.370"
$ws.Range("E11").Value = "  +5.55%  "
$ws.Range("D12").Value = "3.542.61"
$ws.Range("E12").Value = "  +2.76%  "
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("E15").Value = "  +3.96%  "
$ws.Range("D16").Value = "56.674.80"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "3.022.22"
$ws.Range("E17").Value = "  +2.78%  "
$ws.Range("D18").Value = "'5.92"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("D19").Value = "'13.10"
$ws.Range("E19").Value = "  +5.47%  "
$ws.Range("E20").Value = "  +4.12%  "
$ws.Range("D21").Value = "'333.97"
$ws.Range("E21").Value = "  +5.64%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'0.501"
$ws.Range("E23").Value = "  +3.79%  "
$ws.Range("D24").Value = "'64.83"
$ws.Range("E24").Value = "  +4.04%  "
$ws.Range("D25").Value = "3.152.25"
$ws.Range("E25").Value = "  +2.93%  "
$ws.Range("D26").Value = "'0.167"
$ws.Range("E26").Value = "  +4.05%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").Value = "0.0₃0923"
$ws.Range("E28").Value = "  +8.40%  "
$ws.Range("D29").Value = "'6.37"
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("D30").Value = "'6.78"
$ws.Range("E30").Value = "  -3.13%  "
$ws.Range("E31").Value = "  +3.35%  "
$ws.Range("D32").Value = "'20.44"
$ws.Range("E32").Value = "  +2.65%  "
$ws.Range("E33").Value = "  +2.87%  "
$ws.Range("D34").Value = "'153.04"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("B35").Value = "EnergySwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D35").Value = "'27.33"
$ws.Range("E35").Value = "  +15.71%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'4.49"
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("E37").Value = "  +2.76%  "
$ws.Range("E38").Value = "  +2.14%  "
$ws.Range("D39").Value = "'0.0664"
$ws.Range("E39").Value = "  +1.97%  "
$ws.Range("D40").Value = "3.060.24"
$ws.Range("E40").Value = "  +3.06%  "
$ws.Range("D41").Value = "'36.52"
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").Value = "'3.81"
$ws.Range("E43").Value = "  +3.97%  "
$ws.Range("D44").Value = "'0.657"
$ws.Range("E44").Value = "  +3.97%  "
$ws.Range("D45").Value = "2.209.38"
$ws.Range("E45").Value = "  +3.59%  "
$ws.Range("E46").Value = "  +1.49%  "
$ws.Range("E47").Value = "  +6.27%  "
$ws.Range("D48").Value = "'0.933"
$ws.Range("E48").Value = "  +2.65%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'5.84"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'19.74"
$ws.Range("E50").Value = "  +4.65%  "
$ws.Range("E51").Value = "  +1.08%  "
